$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Converted Data")

# --- Row 7 ("Weights"): travel_limit (S7) weight dropped to 0, so the total
#     weight (U7, previously the sum of all policy weights) drops from 13 to 12.
$ws.Range("S7").Value2 = 0
$ws.Range("U7").Value2 = 12

# Reference cell whose formatting (bold / bordered / centered header style,
# used throughout column A) we reuse for every new date row added below.
$srcA = $ws.Range("A221")

# --- Recompute the per-policy "share of active weight" column (U) for every
#     existing data row now that the weight total is 12 instead of 13. The
#     numerator (how many active-weighted policies apply that day) is
#     unchanged, only the denominator changed, so values fall into a handful
#     of contiguous bands.
foreach ($r in 24..25)  { $ws.Range("U$r").Value2 = 0.08333333333333333 }
foreach ($r in 26..26)  { $ws.Range("U$r").Value2 = 0.1666666666666667 }
foreach ($r in 27..34)  { $ws.Range("U$r").Value2 = 0.25 }
foreach ($r in 35..70)  { $ws.Range("U$r").Value2 = 0.5 }
foreach ($r in 71..96)  { $ws.Range("U$r").Value2 = 0.25 }
foreach ($r in 97..221) { $ws.Range("U$r").Value2 = 0.1666666666666667 }

# --- Append 12 new daily rows (9/30/2020 - 10/11/2020), each following the
#     same policy pattern already in place for the most recent days.
$newDates = @("9/30/2020","10/1/2020","10/2/2020","10/3/2020","10/4/2020","10/5/2020","10/6/2020","10/7/2020","10/8/2020","10/9/2020","10/10/2020","10/11/2020")
$rowVals = @{ B=0; C=0; D=1; E=0; F=1; G=0; H=0; I=0; J=0; K=0; L=1; M=0; N=1; O=0; P=0; Q=0; R=0; S=0; T=0 }
$startRow = 222

for ($i = 0; $i -lt $newDates.Count; $i++) {
    $r = $startRow + $i

    # Column A holds the date as literal text (matching the rest of the
    # column), not an auto-converted Excel date serial: force Text format
    # before assigning, then restore the header-cell formatting (border,
    # bold, centered) by copying it over from the row above.
    $dst = $ws.Range("A$r")
    $dst.NumberFormat = "@"
    $dst.Value = $newDates[$i]
    $srcA.Copy()
    $dst.PasteSpecial(-4122)

    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$r").Value2 = $rowVals[$col]
    }
    $ws.Range("U$r").Value2 = 0.1666666666666667
}

$excel.CutCopyMode = 0
